$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (G2)
$wsOverview.Range("G2").Value = "2016-09-06 23:20:08"

# zh-cn sheet: "Correspond Handoff Datetime" (H2) and "Correspond Handback DateTime" (K2)
$wsZhCn.Range("H2").Value = "2016-09-06 23:19:57"
$wsZhCn.Range("K2").Value = "2016-09-06 23:20:33"

# de-de sheet: "Correspond Handoff Datetime" (H2) shares the same underlying
# value as Overview!G2 ("Latest HO Xliff Generate Date"), and
# "Correspond Handback DateTime" (K2)
$wsDeDe.Range("H2").Value = "2016-09-06 23:20:08"
$wsDeDe.Range("K2").Value = "2016-09-06 23:20:41"
